$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.723.46'

$ws.Cells.Item(3, 4).Value = '1.847.64'
$ws.Cells.Item(3, 5).Value = '  +0.15%  '

$ws.Cells.Item(4, 5).Value = '  +0.06%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '313.81'
$ws.Cells.Item(5, 5).Value = '  -0.41%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.002'
$ws.Cells.Item(6, 5).Value = '  +0.11%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4320'
$ws.Cells.Item(7, 5).Value = '  +1.16%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3655'
$ws.Cells.Item(8, 5).Value = '  +0.20%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.07344'
$ws.Cells.Item(9, 5).Value = '  +0.86%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.8801'
$ws.Cells.Item(10, 5).Value = '  -1.70%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '20.82'
$ws.Cells.Item(11, 5).Value = '  +0.65%  '

$ws.Cells.Item(12, 4).Value = '1.874.39'
$ws.Cells.Item(12, 5).Value = '  +0.21%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '5.363'
$ws.Cells.Item(13, 5).Value = '  -0.51%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.543'
$ws.Cells.Item(14, 5).Value = '  -0.49%  '

$ws.Cells.Item(15, 5).Value = '  +0.82%  '

$ws.Cells.Item(16, 5).Value = '  +0.10%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '79.91'
$ws.Cells.Item(17, 5).Value = '  +1.83%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.000009020'
$ws.Cells.Item(18, 5).Value = '  +1.90%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '1.002'
$ws.Cells.Item(19, 5).Value = '  +0.03%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '15.41'
$ws.Cells.Item(20, 5).Value = '  -0.91%  '

$ws.Cells.Item(21, 4).Value = '27.753.88'
$ws.Cells.Item(21, 5).Value = '  +0.46%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '4.985'
$ws.Cells.Item(22, 5).Value = '  +0.01%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '10.35'
$ws.Cells.Item(23, 5).Value = '  -1.83%  '

$ws.Cells.Item(24, 4).Value = '2.085.76'
$ws.Cells.Item(24, 5).Value = '  +0.66%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '1.996'
$ws.Cells.Item(25, 5).Value = '  -2.20%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '155.79'
$ws.Cells.Item(26, 5).Value = '  +0.52%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '18.66'
$ws.Cells.Item(27, 5).Value = '  +1.29%  '

$ws.Cells.Item(28, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '5.269'
$ws.Cells.Item(28, 5).Value = '  +0.48%  '

$ws.Cells.Item(29, 2).Value = 'BitcoinCash'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '119.97'
$ws.Cells.Item(29, 5).Value = '  +4.91%  '

$ws.Cells.Item(30, 5).Value = '  +2.63%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.08896'
$ws.Cells.Item(31, 5).Value = '  -0.02%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.7604'
$ws.Cells.Item(32, 5).Value = '  -2.97%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '4.559'
$ws.Cells.Item(33, 5).Value = '  -0.53%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '2.962'
$ws.Cells.Item(34, 5).Value = '  -0.36%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.129'
$ws.Cells.Item(35, 5).Value = '  +2.12%  '

$ws.Cells.Item(36, 5).Value = '  +0.05%  '

$ws.Cells.Item(37, 2).Value = 'Hedera'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.05451'
$ws.Cells.Item(37, 5).Value = '  +0.14%  '

$ws.Cells.Item(38, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.110'
$ws.Cells.Item(38, 5).Value = '  +0.77%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.01940'
$ws.Cells.Item(39, 5).Value = '  +0.56%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.839'
$ws.Cells.Item(40, 5).Value = '  +2.06%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.5103'
$ws.Cells.Item(41, 5).Value = '  +0.58%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.1670'
$ws.Cells.Item(42, 5).Value = '  +1.19%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '6.641'
$ws.Cells.Item(43, 5).Value = '  -2.71%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '8.421'
$ws.Cells.Item(44, 5).Value = '  +2.12%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.06555'
$ws.Cells.Item(45, 5).Value = '  -1.30%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '10.36'
$ws.Cells.Item(46, 5).Value = '  -0.23%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '105.71'
$ws.Cells.Item(47, 5).Value = '  -0.05%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.4670'
$ws.Cells.Item(48, 5).Value = '  -0.99%  '

$ws.Cells.Item(49, 5).Value = '  +0.12%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.643'
$ws.Cells.Item(50, 5).Value = '  +0.46%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '64.67'
$ws.Cells.Item(51, 5).Value = '  +0.25%  '
